# Covid_19_Dataset_and_References/References/34.xlsx - "Going through the dataset, updating"
# Adds a new "Other found locations" column (I) and refreshes the ID lookup
# for row 2 (not found) and the Authors lists for rows 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: "Other found locations" ---
$ws.Range("I1").Value = "Other found locations"
$ws.Range("I2").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("I4").Value = "_PMC_elsevier"
$ws.Range("I5").Value = "_PMC"

# --- Row 2: ID lookup now comes back empty ---
$ws.Range("F2").Value = "not found"
$ws.Range("G2").Value = "N/A"

# --- Rows 4 & 5: refreshed Authors strings ---
$ws.Range("E4").Value = '[Ruchong%Chen%NULL%0,   Wenhua%Liang%NULL%4,   Mei%Jiang%NULL%3,   Weijie%Guan%NULL%4,   Chen%Zhan%NULL%1,   Tao%Wang%NULL%0,   Chunli%Tang%NULL%2,   Ling%Sang%NULL%0,   Jiaxing%Liu%NULL%1,   Zhengyi%Ni%NULL%3,   Yu%Hu%NULL%0,   Lei%Liu%NULL%0,   Hong%Shan%NULL%0,   Chunliang%Lei%NULL%3,   Yixiang%Peng%NULL%3,   Li%Wei%NULL%0,   Yong%Liu%NULL%0,   Yahua%Hu%NULL%3,   Peng%Peng%NULL%0,   Jianming%Wang%NULL%0,   Jiyang%Liu%NULL%3,   Zhong%Chen%NULL%0,   Gang%Li%NULL%0,   Zhijian%Zheng%NULL%3,   Shaoqin%Qiu%NULL%3,   Jie%Luo%NULL%0,   Changjiang%Ye%NULL%3,   Shaoyong%Zhu%NULL%3,   Xiaoqing%Liu%NULL%1,   Linling%Cheng%NULL%1,   Feng%Ye%NULL%0,   Jinping%Zheng%NULL%3,   Nuofu%Zhang%NULL%3,   Yimin%Li%NULL%3,   Jianxing%He%NULL%3,   Shiyue%Li%lishiyue@188.com%0,   Nanshan%Zhong%NULL%5,   NULL%NULL%NULL%0]'
$ws.Range("E5").Value = '[Wei-jie%Guan%NULL%0,   Zheng-yi%Ni%NULL%0,   Zheng-yi%Ni%NULL%0,   Yu%Hu%NULL%0,   Wen-hua%Liang%NULL%0,   Chun-quan%Ou%NULL%0,   Jian-xing%He%NULL%0,   Lei%Liu%NULL%0,   Hong%Shan%NULL%0,   Chun-liang%Lei%NULL%0,   David S.C.%Hui%NULL%0,   Bin%Du%NULL%0,   Lan-juan%Li%NULL%0,   Guang%Zeng%NULL%0,   Kwok-Yung%Yuen%NULL%0,   Ru-chong%Chen%NULL%0,   Chun-li%Tang%NULL%0,   Tao%Wang%NULL%0,   Ping-yan%Chen%NULL%0,   Jie%Xiang%NULL%0,   Shi-yue%Li%NULL%0,   Jin-lin%Wang%NULL%0,   Zi-jing%Liang%NULL%0,   Yi-xiang%Peng%NULL%0,   Li%Wei%NULL%0,   Yong%Liu%NULL%0,   Ya-hua%Hu%NULL%0,   Peng%Peng%NULL%0,   Jian-ming%Wang%NULL%0,   Ji-yang%Liu%NULL%0,   Zhong%Chen%NULL%0,   Gang%Li%NULL%0,   Zhi-jian%Zheng%NULL%0,   Shao-qin%Qiu%NULL%0,   Jie%Luo%NULL%0,   Chang-jiang%Ye%NULL%0,   Shao-yong%Zhu%NULL%0,   Nan-shan%Zhong%NULL%0]'
